# Appends trade-sampling rows (rows 5-20) to the results sheet, matching
# the "pick trade every 1/2/4 hours" recording logic.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @(4, 16.88810595718799, 7.410636442894507, 1.595419847328244, 4588),
    @(5, -23.18958569898614, 12.51362546326575, 1.174242424242424, 4587),
    @(6, 63.13150677758502, 9.136502398604449, 1.479289940828402, 4586),
    @(7, 61.20206252164174, 8.287895310796074, 1.620689655172414, 4585),
    @(8, 11.47957239204662, 3.926701570680628, 1.535211267605634, 4584),
    @(9, 54.95367048012035, 6.938686449923631, 1.355555555555556, 4583),
    @(10, 102.3194526587594, 6.678306416412047, 1.615384615384615, 4582),
    @(11, 74.99318399046976, 6.090373280943026, 1.607476635514019, 4581),
    @(12, 84.77249498106892, 5.349344978165939, 1.692307692307692, 4580),
    @(13, 58.50058886331681, 2.839047827036471, 1.708333333333333, 4579),
    @(14, 50.44483745362236, 4.084753167321974, 1.367088607594937, 4578),
    @(15, 68.70088120721429, 4.610006554511689, 1.573170731707317, 4577),
    @(16, 7.931818834894001, 4.567307692307692, 1.177083333333333, 4576),
    @(17, 83.65057414842337, 3.497267759562841, 1.388059701492537, 4575),
    @(18, 16.64545044202842, 3.279405334499344, 1.238805970149254, 4574),
    @(19, -26.39520180019328, 4.176689263065821, 1.122222222222222, 4573)
)

$startRow = 5
for ($i = 0; $i -lt $data.Count; $i++) {
    $r = $startRow + $i
    $row = $data[$i]
    $ws.Cells.Item($r, 1).Value = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
    $ws.Cells.Item($r, 4).Value = $row[3]
    $ws.Cells.Item($r, 5).Value = $row[4]
}
